# "Trafikkort - add on startet."
# Applies:
#  1. Rename the "Google Analytics - Trafikkort" output-file setting value
#     to "Google Analytics - Trafikinfo" on the "Indstillinger" sheet.
#  2. Add a new worksheet "Data modifikation 2" (after "Data modifikation")
#     with threshold / summarize-and-delete (event category) settings, and
#     make it the active sheet/tab.

$wb = $excel.ActiveWorkbook

# --- 1. Update the "Navn på output excel fil" value -----------------------
$settings = $wb.Worksheets.Item("Indstillinger")
$settings.Range("B6").Value = "Google Analytics - Trafikinfo"

# --- 2. Add the new "Data modifikation 2" worksheet ------------------------
$dataMod = $wb.Worksheets.Item("Data modifikation")
$newSheet = $wb.Worksheets.Add($Null, $dataMod)
$newSheet.Name = "Data modifikation 2"

$newSheet.Range("A1").Value = "Summér og slet (hændelseskategori)"
$newSheet.Range("A1").Font.Bold = $true

$newSheet.Range("A2").Value = "Variabel slut navn (bibeholdes)"
$newSheet.Range("B2").Value = "Variabel der starter med:"

$newSheet.Range("A3").Value = "https://trafikkort.vejdirektoratet.dk/?utm_source=danwest.de&utm_campaign=3ec1f9fb03-EMAIL_CAMPAIGN_2018_06_11_08"
$newSheet.Range("B3").Value = "https://trafikkort.vejdirektoratet.dk/?utm_source=danwest.de&utm_campaign=3ec1f9fb03-EMAIL_CAMPAIGN_2018_06_11_08"

$newSheet.Range("E1").Value = "Threshold"
$newSheet.Range("E2").Value = "Opret separat ark for værdier over"

# NOTE: Excel (and this COM runtime) stores column widths quantized to
# 1/6-character (pixel) steps, so the inputs below are chosen to land on
# the quantized value nearest the target widths (116.140625 / 15.140625 /
# 31.85546875 "characters") actually present in the target workbook.
$newSheet.Columns.Item(1).ColumnWidth = 115.33
$newSheet.Columns.Item(2).ColumnWidth = 14.33
$newSheet.Columns.Item(5).ColumnWidth = 31.00

$newSheet.PageSetup.PaperSize = 9
$newSheet.PageSetup.Orientation = 1

# Selections matching the target workbook
$settings.Range("B42").Select() | Out-Null
$dataMod.Range("E43").Select() | Out-Null
$newSheet.Range("A5").Select() | Out-Null

# Make the new sheet the active / visible tab
$newSheet.Activate()
